$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.004.13"
$ws.Range("D3").Value = "3.527.15"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.01"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.19"
$ws.Range("E6").Value = "  +4.89%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.528.88"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("E10").Value = "  +5.93%  "
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "4.139.76"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.48"
$ws.Range("E14").Value = "  +11.04%  "
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "67.964.79"
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "3.536.00"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.55"
$ws.Range("E20").Value = "  +2.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "401.58"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.04"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.12"
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.547"
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.56"
$ws.Range("E28").Value = "  +2.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.179"
$ws.Range("E29").Value = "  -2.47%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.32"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.05"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.58"
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.44"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.885"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.81"
$ws.Range("E41").Value = "  +7.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.02"
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.74"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").Value = "2.899.97"
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.65"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0738"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.96"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.64"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "352.73"
$ws.Range("E49").Value = "  +3.68%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("E51").Value = "  -1.32%  "
